# "chore: adapt column header formatting to respective input file names"
#
# 1) Rename the header row:
#      *_old -> *_FV2404
#      *_new -> *_FV2410
#    (the "diff" header is left untouched)
# 2) Freeze the header row.
# 3) Turn the data range into a proper Excel Table ("Table1") whose
#    column names mirror the (renamed) header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerRange = $ws.Range("A1:U1")
$colCount = $headerRange.Columns.Count

for ($c = 1; $c -le $colCount; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value()
    if ($v -like "*_old") {
        $cell.Value = ($v -replace "_old$", "_FV2404")
    } elseif ($v -like "*_new") {
        $cell.Value = ($v -replace "_new$", "_FV2410")
    }
}

# Freeze the header row (row 1) in the sheet view.
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into a native Excel Table ("Table1") so the
# (now renamed) header row drives the table column names / autofilter.
$dataRange = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"
